$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 111633890
$ws.Range("B4").Value = 90658
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 4361
$ws.Range("F4").Value = "Orange taggsvamp"
$ws.Range("G4").Value = "Hydnellum aurantiacum"
$ws.Range("H4").Value = "(Batsch:Fr.) P.Karst."

$ws.Range("A5").Value = 111634304
$ws.Range("B5").Value = 90687
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 5964
$ws.Range("F5").Value = "Fjällig taggsvamp s.str."
$ws.Range("G5").Value = "Sarcodon imbricatus s.str."
$ws.Range("H5").Value = "(L.:Fr.) P.Karst."
$ws.Range("Q5").Value = 676708.8668162767
$ws.Range("R5").Value = 6618511.450801066

$ws.Range("A6").Value = 111633843
$ws.Range("B6").Value = 90687
$ws.Range("E6").Value = 5964
$ws.Range("F6").Value = "Fjällig taggsvamp s.str."
$ws.Range("G6").Value = "Sarcodon imbricatus s.str."
$ws.Range("H6").Value = "(L.:Fr.) P.Karst."

$ws.Range("A8").Value = 111633837
$ws.Range("B8").Value = 98535
$ws.Range("E8").Value = 222498
$ws.Range("F8").Value = "Blåsippa"
$ws.Range("G8").Value = "Hepatica nobilis"
$ws.Range("H8").Value = "Schreb."
$ws.Range("Q8").Value = 676486.710397501
$ws.Range("R8").Value = 6618439.724061669
